# Updated TestData for Portugal Market
#
# - Germany sheet: selection moved to A12 (cosmetic, matches authored diff)
# - New "Portugal" worksheet added after "Swiss", cloned from "Swiss" and
#   then re-pointed to Portugal-specific market code / JIRA ticket, with the
#   "PROFILE Communicator" row moved to the bottom of its list and highlighted
#   with a black-font / boxed-border style.

$wb = $excel.ActiveWorkbook

# --- Germany: just a cosmetic selection move (A10 -> A12) -------------------
$germany = $wb.Worksheets.Item("Germany")
[void]$germany.Activate()
[void]$germany.Range("A12").Select()

# --- Clone "Swiss" into a new "Portugal" sheet, placed right after it -------
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy([System.Reflection.Missing]::Value, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Market name + JIRA/test-case reference for Portugal
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2404"

# Re-order the printer list: "PROFILE Communicator" drops to the end of the
# middle block (row 14) so the rest shift up one slot.
$portugal.Range("A11").Value = "MOXA Event Logger"
$portugal.Range("A12").Value = "PX-PR"
$portugal.Range("A13").Value = "RS800"
$portugal.Range("A14").Value = "PROFILE Communicator"

# Highlight the relocated "PROFILE Communicator" row: explicit black font and
# a thin box border on the right/top/bottom (no left edge).
$a14 = $portugal.Range("A14")
$a14.Font.Color = 0
$a14.Borders.Item(7).LineStyle = -4142
$a14.Borders.Item(10).LineStyle = 1
$a14.Borders.Item(8).LineStyle = 1
$a14.Borders.Item(9).LineStyle = 1

# Make Portugal the active tab with A14 selected (last, so it "wins" as the
# workbook's active sheet / activeTab).
[void]$portugal.Activate()
[void]$a14.Select()
